# Updated cryptos list on Mon Aug  5 15:42:17 UTC 2024 with GitHub Actions
#
# Applies the refreshed price/volume(1h) figures scraped for this run.
# Column D ("Price") holds values that look numeric (e.g. "1.00", "5.70")
# but are stored as plain text in the source data (note multi-dot values
# like "54.457.69" which are not valid numbers at all). Excel's COM layer
# auto-coerces clean numeric-looking strings assigned via .Value into real
# numbers, which would silently mangle figures like "1.00" -> 1. To keep
# every Price cell as text (matching the original workbook), the D2:D51
# range is temporarily switched to a Text number format before the writes,
# then restored to the default "Normal" style afterwards so no stray
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '54.457.69'
$ws.Range("E2").Value = '  -8.11%  '
$ws.Range("D3").Value = '2.419.01'
$ws.Range("E3").Value = '  -14.75%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '465.31'
$ws.Range("E5").Value = '  -8.01%  '
$ws.Range("D6").Value = '130.50'
$ws.Range("E6").Value = '  -3.85%  '
$ws.Range("E7").Value = '  -0.67%  '
$ws.Range("D8").Value = '0.489'
$ws.Range("E8").Value = '  -7.32%  '
$ws.Range("D9").Value = '2.428.75'
$ws.Range("E9").Value = '  -14.24%  '
$ws.Range("D10").Value = '0.0947'
$ws.Range("E10").Value = '  -8.59%  '
$ws.Range("D11").Value = '5.37'
$ws.Range("E11").Value = '  -9.38%  '
$ws.Range("D12").Value = '0.318'
$ws.Range("E12").Value = '  -8.73%  '
$ws.Range("E13").Value = '  -4.13%  '
$ws.Range("D14").Value = '2.824.62'
$ws.Range("E14").Value = '  -15.08%  '
$ws.Range("D15").Value = '54.232.06'
$ws.Range("E15").Value = '  -8.63%  '
$ws.Range("D16").Value = '19.71'
$ws.Range("E16").Value = '  -9.70%  '
$ws.Range("D17").Value = '0.0000130'
$ws.Range("E17").Value = '  -4.02%  '
$ws.Range("D18").Value = '2.416.56'
$ws.Range("E18").Value = '  -14.91%  '
$ws.Range("D19").Value = '4.22'
$ws.Range("E19").Value = '  -11.11%  '
$ws.Range("D20").Value = '312.14'
$ws.Range("E20").Value = '  -11.31%  '
$ws.Range("D21").Value = '9.40'
$ws.Range("E21").Value = '  -15.15%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("D23").Value = '5.70'
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("D24").Value = '5.41'
$ws.Range("E24").Value = '  -13.39%  '
$ws.Range("D25").Value = '56.73'
$ws.Range("E25").Value = '  -10.08%  '
$ws.Range("D26").Value = '1.01'
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("D27").Value = '0.156'
$ws.Range("E27").Value = '  -9.50%  '
$ws.Range("D28").Value = '0.382'
$ws.Range("E28").Value = '  -11.05%  '
$ws.Range("D29").Value = '2.496.65'
$ws.Range("E29").Value = '  -15.90%  '
$ws.Range("D30").Value = '7.19'
$ws.Range("E30").Value = '  -3.28%  '
$ws.Range("D31").Value = '0.995'
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").Value = '0.0₃0710'
$ws.Range("E32").Value = '  -12.75%  '
$ws.Range("D33").Value = '150.34'
$ws.Range("E33").Value = '  -0.33%  '
$ws.Range("D34").Value = '17.69'
$ws.Range("E34").Value = '  -6.95%  '
$ws.Range("D35").Value = '1.41'
$ws.Range("E35").Value = '  -12.83%  '
$ws.Range("D36").Value = '5.04'
$ws.Range("E36").Value = '  -6.10%  '
$ws.Range("D37").Value = '3.54'
$ws.Range("E37").Value = '  -15.26%  '
$ws.Range("D38").Value = '1.05'
$ws.Range("E38").Value = '  -7.57%  '
$ws.Range("D39").Value = '0.804'
$ws.Range("E39").Value = '  -13.04%  '
$ws.Range("D40").Value = '33.63'
$ws.Range("E40").Value = '  -8.05%  '
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("D42").Value = '0.599'
$ws.Range("E42").Value = '  -4.88%  '
$ws.Range("D43").Value = '0.0531'
$ws.Range("E43").Value = '  -5.29%  '
$ws.Range("D44").Value = '3.30'
$ws.Range("E44").Value = '  -6.50%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = '1.25'
$ws.Range("E45").Value = '  -8.91%  '
$ws.Range("B46").Value = 'WhiteBITCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D46").Value = '10.14'
$ws.Range("E46").Value = '  -1.93%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.970.34'
$ws.Range("E47").Value = '  -11.35%  '
$ws.Range("D48").Value = '0.0220'
$ws.Range("E48").Value = '  -3.03%  '
$ws.Range("D49").Value = '0.0869'
$ws.Range("E49").Value = '  -2.04%  '
$ws.Range("D50").Value = '4.31'
$ws.Range("E50").Value = '  -4.51%  '
$ws.Range("D51").Value = '16.51'
$ws.Range("E51").Value = '  -15.27%  '

$priceRange.Style = "Normal"
